$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.00041467678966000676
$ws.Range("A3").Value = 0.00023120432160794735
$ws.Range("H3").Value = 5.5
$ws.Range("A4").Value = 0.00018347245350014418
$ws.Range("H4").Value = 5.390532970428467
